$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data row (row 37): date 2024-11-25 (serial 45621), and score updates
$ws.Range("A37").Value = 45621
$ws.Range("B37").Value = 94
$ws.Range("C37").Value = 76
$ws.Range("D37").Value = 89

# Match number formatting: row 37's date cell (A37) takes on the "last row" style
# previously held by A36, while A36 moves to the standard date style.
$ws.Range("A36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A37").NumberFormat = "YYYY-MM-DD"
